# Apply the commit: insert two new price-report rows for Tuna right after
# row 331 (i.e. they become the new rows 332 and 333), pushing all the
# existing rows 332-440 down to rows 334-442.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 332 (each insert pushes everything
# below down by one row).
$ws.Rows.Item(332).Insert()
$ws.Rows.Item(332).Insert()

# ---- New row 332 ----
$ws.Range("A332").Value = 6
$ws.Range("B332").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C332").Value = "Metropolitana"
$ws.Range("D332").Value = 44809
$ws.Range("E332").Value = 13
$ws.Range("F332").Value = "Fruta"
$ws.Range("G332").Value = 100107
$ws.Range("H332").Value = "Otros"
$ws.Range("I332").Value = 100107011
$ws.Range("J332").Value = "Tuna"
$ws.Range("K332").Value = "Sin especificar"
$ws.Range("L332").Value = "Extra (doble especial)"
$ws.Range("M332").Value = 200
$ws.Range("N332").Value = 28000
$ws.Range("O332").Value = 28000
$ws.Range("P332").Value = 28000
$ws.Range("Q332").Value = "`$/caja 18 kilos granel"
$ws.Range("R332").Value = "Provincia de Melipilla"
$ws.Range("S332").Value = 1556
$ws.Range("T332").Value = 18

# ---- New row 333 ----
$ws.Range("A333").Value = 6
$ws.Range("B333").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C333").Value = "Metropolitana"
$ws.Range("D333").Value = 44809
$ws.Range("E333").Value = 13
$ws.Range("F333").Value = "Fruta"
$ws.Range("G333").Value = 100107
$ws.Range("H333").Value = "Otros"
$ws.Range("I333").Value = 100107011
$ws.Range("J333").Value = "Tuna"
$ws.Range("K333").Value = "Sin especificar"
$ws.Range("L333").Value = "Primera"
$ws.Range("M333").Value = 40
$ws.Range("N333").Value = 24000
$ws.Range("O333").Value = 24000
$ws.Range("P333").Value = 24000
$ws.Range("Q333").Value = "`$/caja 18 kilos granel"
$ws.Range("R333").Value = "Provincia de Melipilla"
$ws.Range("S333").Value = 1333
$ws.Range("T333").Value = 18

# Make sure the date cells keep the date number format used by the rest
# of column D (style index 2 in the original workbook), same as the
# cells immediately below, which the Insert() already applied; this is
# a safety net in case it wasn't carried over.
$ws.Range("D332").NumberFormat = $ws.Range("D334").NumberFormat
$ws.Range("D333").NumberFormat = $ws.Range("D334").NumberFormat
